$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Copy the date formatting from an existing date cell so the new date
# cells pick up the same (already-existing) style instead of creating a
# brand new style entry.
$ws.Range("C2").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 13: The Pioneers by David McCullough
$ws.Range("A13").Value = "The Pioneers"
$ws.Range("B13").Value = "David McCullough"
$ws.Range("C13").Value2 = 43847
$ws.Range("D13").Value2 = 43850
$ws.Range("E13").Value = "history;ohio;pioneers;america;1700s;1800s"
$ws.Range("F13").Value = "Audio"
$ws.Range("G13").Value = "10 Hrs 33 Mins"

# Row 14: Iaccoca: An Autobiography by Lee Iaccoca (in progress - only
# a start date so far, no finish date or length yet)
$ws.Range("A14").Value = "Iaccoca: An Autobiography"
$ws.Range("B14").Value = "Lee Iaccoca"
$ws.Range("C14").Value2 = 43847
$ws.Range("E14").Value = "memoir;business;ford;chrysler;Lee Iaccoca"
$ws.Range("F14").Value = "Hard Copy"

$ws.Range("A15").Select()
